$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 06:42"

# India (row 5) gets updated figures
$ws.Cells.Item(5, 2).Value = 7550273
$ws.Cells.Item(5, 3).Value = 2035
$ws.Cells.Item(5, 4).Value = 6663608
$ws.Cells.Item(5, 5).Value = 772023

# Honduras overtakes Bielorrusia in the ranking, so the two rows swap places.
# Row 53 becomes Honduras (with its new, larger figures)
$ws.Cells.Item(53, 1).Value = "Honduras"
$ws.Cells.Item(53, 2).Value = 88425
$ws.Cells.Item(53, 3).Value = 831
$ws.Cells.Item(53, 4).Value = 34964
$ws.Cells.Item(53, 5).Value = 50893
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 5
$ws.Cells.Item(53, 8).Value = 2568

# Row 54 becomes Bielorrusia (its previous, unchanged figures)
$ws.Cells.Item(54, 1).Value = "Bielorrusia"
$ws.Cells.Item(54, 2).Value = 87698
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 79757
$ws.Cells.Item(54, 5).Value = 7012
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 929

# Tailandia (row 147) updated figures
$ws.Cells.Item(147, 2).Value = 3691
$ws.Cells.Item(147, 3).Value = 5
$ws.Cells.Item(147, 4).Value = 3488
$ws.Cells.Item(147, 5).Value = 144

# Butan (row 186) updated figures
$ws.Cells.Item(186, 2).Value = 327
$ws.Cells.Item(186, 3).Value = 2
$ws.Cells.Item(186, 4).Value = 301
